# Daily attendance processing - 2025-12-27 15:54:12
# For every data row in the "Recorded By" column (G), swap the first two
# comma-separated entries (e.g. "System, user@example.com" becomes
# "user@example.com, System"), leaving any remaining entries in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $newVal = [string]::Join(", ", $parts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
